$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.719.89"
$ws.Range("E2").Value = "  +0.64%  "
$ws.Range("D3").Value = "1.950.88"
$ws.Range("E3").Value = "  +1.89%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "247.40"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.39%  "
$ws.Range("E6").Value = "  -0.04%  "
$ws.Range("E7").Value = "  -0.59%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2949"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.97%  "
$ws.Range("E9").Value = "  +1.68%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "112.23"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.23%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "19.50"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.49%  "
$ws.Range("D12").Value = "1.944.99"
$ws.Range("E12").Value = "  +1.53%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.538"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +5.33%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.07647"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.79%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6931"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.27%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "295.99"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +9.59%  "
$ws.Range("D17").Value = "30.704.36"
$ws.Range("E17").Value = "  +0.64%  "
$ws.Range("E18").Value = "  +4.05%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.706"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.95%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007699"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.27%  "
$ws.Range("D21").Value = "2.193.88"
$ws.Range("E21").Value = "  +1.26%  "
$ws.Range("E22").Value = "  -0.01%  "
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.561"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.52%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.792"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.27%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "168.18"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.75%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.37"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.77%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.181"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.14%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.1093"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.24%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.438"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.59%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.748"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +17.75%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.445"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +7.97%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05083"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.07%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7797"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +7.18%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.164"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.60%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02073"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.37%  "
$ws.Range("E37").Value = "  +0.72%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.705"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.30%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.043"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.83%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "110.88"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.16%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4458"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.13%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8746"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.11%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.948"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.57%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "70.99"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.98%  "
$ws.Range("E45").Value = "  +0.18%  "
$ws.Range("E46").Value = "  +0.52%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.514"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.60%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "48.73"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.15%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.1252"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.76%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "35.50"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.84%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.2551"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.28%  "
